$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric values in columns B:E (rows 2-13) to the nearest integer,
# matching the commit's change to write Pot_/Ontpl_ data as integers.
$range = $ws.Range("B2:E13")
foreach ($cell in $range.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = [Math]::Round([double]$val)
    }
}
